$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N ("Late"), pushing
# Late/Outstanding(heading)/Outstanding one column to the right.
$ws.Columns("N").Insert()

# Select the new cell and make this the active sheet/tab, matching the
# saved view state after the edit.
$ws.Select()
$ws.Range("L15").Select()
